$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5999.1113
$ws.Range("J86").Value = 5999.25
$ws.Range("L86").Value = 5999.25
$ws.Range("N86").Value = -8245.25
$ws.Range("H89").Value = 5999.1113
$ws.Range("J89").Value = 5999.25
$ws.Range("L89").Value = 29996.25
$ws.Range("N89").Value = -41228.25
$ws.Range("H111").Value = 2485
$ws.Range("I111").Value = 2535.7778
$ws.Range("J111").Value = 2332.6667
$ws.Range("K111").Value = 7607.3334
$ws.Range("L111").Value = 6998.000100000001
$ws.Range("M111").Value = -4540.3334
$ws.Range("N111").Value = -13132.0001
$ws.Range("H132").Value = 2258.2
$ws.Range("I132").Value = 2219.158
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 6657.474
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -4127.474
$ws.Range("N132").Value = -14060
$ws.Range("H135").Value = 2226.375
$ws.Range("I135").Value = 1547.25
$ws.Range("J135").Value = 2905.5
$ws.Range("K135").Value = 13925.25
$ws.Range("L135").Value = 26149.5
$ws.Range("M135").Value = -11390.25
$ws.Range("N135").Value = -31219.5
$ws.Range("H138").Value = 4013.7083
$ws.Range("I138").Value = 1322.6666
$ws.Range("K138").Value = 3967.9998
$ws.Range("M138").Value = 1172.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8155.684
$ws.Range("I32").Value = 3497.9666
$ws.Range("J32").Value = 25622.125
$ws.Range("K32").Value = 3497.9666
$ws.Range("L32").Value = 25622.125
$ws.Range("M32").Value = -3210.9666
$ws.Range("N32").Value = -26196.125
$ws.Range("H46").Value = 25020.6
$ws.Range("I46").Value = 26276
$ws.Range("J46").Value = 19999
$ws.Range("K46").Value = 26276
$ws.Range("L46").Value = 19999
$ws.Range("M46").Value = -25957
$ws.Range("N46").Value = -20637
$ws.Range("H61").Value = 4654.5
$ws.Range("I61").Value = 4604.6665
$ws.Range("J61").Value = 4704.3335
$ws.Range("K61").Value = 4604.6665
$ws.Range("L61").Value = 4704.3335
$ws.Range("M61").Value = -4392.6665
$ws.Range("N61").Value = -5128.3335
$ws.Range("H74").Value = 1605.9836
$ws.Range("I74").Value = 1210.585
$ws.Range("J74").Value = 4225.5
$ws.Range("K74").Value = 1210.585
$ws.Range("L74").Value = 4225.5
$ws.Range("M74").Value = -336.585
$ws.Range("N74").Value = -5973.5
$ws.Range("H77").Value = 1605.9836
$ws.Range("I77").Value = 1210.585
$ws.Range("J77").Value = 4225.5
$ws.Range("K77").Value = 6052.925
$ws.Range("L77").Value = 21127.5
$ws.Range("M77").Value = -1684.925
$ws.Range("N77").Value = -29863.5
$ws.Range("H122").Value = 3558.5
$ws.Range("I122").Value = 2425.8333
$ws.Range("J122").Value = 4408
$ws.Range("K122").Value = 7277.499899999999
$ws.Range("L122").Value = 13224
$ws.Range("M122").Value = -4827.499899999999
$ws.Range("N122").Value = -18124
$ws.Range("H132").Value = 1574.6666
$ws.Range("I132").Value = 1602.5883
$ws.Range("J132").Value = 1100
$ws.Range("K132").Value = 4807.7649
$ws.Range("L132").Value = 3300
$ws.Range("M132").Value = -2277.7649
$ws.Range("N132").Value = -8360
$ws.Range("H136").Value = 4654.5
$ws.Range("I136").Value = 4604.6665
$ws.Range("J136").Value = 4704.3335
$ws.Range("K136").Value = 13813.9995
$ws.Range("L136").Value = 14113.0005
$ws.Range("M136").Value = -11263.9995
$ws.Range("N136").Value = -19213.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1694.1818
$ws.Range("I86").Value = 1548.1111
$ws.Range("K86").Value = 1548.1111
$ws.Range("M86").Value = -425.1111000000001
$ws.Range("H89").Value = 1694.1818
$ws.Range("I89").Value = 1548.1111
$ws.Range("K89").Value = 7740.5555
$ws.Range("M89").Value = -2124.5555
$ws.Range("H94").Value = 803.26666
$ws.Range("I94").Value = 803.26666
$ws.Range("K94").Value = 803.26666
$ws.Range("M94").Value = -352.26666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 276.45456
$ws.Range("I7").Value = 254.1
$ws.Range("K7").Value = 254.1
$ws.Range("M7").Value = -141.1
$ws.Range("H31").Value = 2947.075
$ws.Range("I31").Value = 1531.5518
$ws.Range("J31").Value = 6678.909
$ws.Range("K31").Value = 1531.5518
$ws.Range("L31").Value = 6678.909
$ws.Range("M31").Value = -1236.5518
$ws.Range("N31").Value = -7268.909
$ws.Range("H34").Value = 2947.075
$ws.Range("I34").Value = 1531.5518
$ws.Range("J34").Value = 6678.909
$ws.Range("K34").Value = 1531.5518
$ws.Range("L34").Value = 6678.909
$ws.Range("M34").Value = -1329.5518
$ws.Range("N34").Value = -7082.909
$ws.Range("H122").Value = 3654.125
$ws.Range("J122").Value = 2755.75
$ws.Range("L122").Value = 8267.25
$ws.Range("N122").Value = -13167.25
$ws.Range("H132").Value = 3604.182
$ws.Range("I132").Value = 3181.4443
$ws.Range("K132").Value = 9544.332900000001
$ws.Range("M132").Value = -7014.332900000001
$ws.Range("H134").Value = 2404.9756
$ws.Range("I134").Value = 1690.4286
$ws.Range("K134").Value = 5071.2858
$ws.Range("M134").Value = -2536.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1514
$ws.Range("I3").Value = 1514
$ws.Range("K3").Value = 4542
$ws.Range("M3").Value = -4430
$ws.Range("H5").Value = 452.33334
$ws.Range("I5").Value = 499.55554
$ws.Range("J5").Value = 424
$ws.Range("K5").Value = 1498.66662
$ws.Range("L5").Value = 1272
$ws.Range("M5").Value = -1386.66662
$ws.Range("N5").Value = -1496
$ws.Range("H11").Value = 583.6
$ws.Range("I11").Value = 167.2
$ws.Range("J11").Value = 1000
$ws.Range("K11").Value = 501.6
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = -361.6
$ws.Range("N11").Value = -3280
$ws.Range("H107").Value = 511.06668
$ws.Range("J107").Value = 504.66666
$ws.Range("L107").Value = 1513.99998
$ws.Range("N107").Value = -5353.999980000001
$ws.Range("H113").Value = 1186.25
$ws.Range("I113").Value = 430.66666
$ws.Range("J113").Value = 1639.6
$ws.Range("K113").Value = 1291.99998
$ws.Range("L113").Value = 4918.799999999999
$ws.Range("N113").Value = -9258.799999999999
$ws.Range("M113").Value = 878.0000199999999
$ws.Range("H122").Value = 440.2857
$ws.Range("I122").Value = 461.66666
$ws.Range("K122").Value = 4154.99994
$ws.Range("M122").Value = -1704.99994
$ws.Range("H132").Value = 2500
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 22500
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -27560
$ws.Range("H134").Value = 9015
$ws.Range("I134").Value = 8030
$ws.Range("K134").Value = 24090
$ws.Range("M134").Value = -19020
$ws.Range("H135").Value = 452.33334
$ws.Range("I135").Value = 499.55554
$ws.Range("J135").Value = 424
$ws.Range("K135").Value = 4495.99986
$ws.Range("L135").Value = 3816
$ws.Range("M135").Value = -1960.99986
$ws.Range("N135").Value = -8886
$ws.Range("H136").Value = 9636.223
$ws.Range("I136").Value = 3818.4285
$ws.Range("J136").Value = 29998.5
$ws.Range("K136").Value = 11455.2855
$ws.Range("L136").Value = 89995.5
$ws.Range("M136").Value = -6355.2855
$ws.Range("N136").Value = -100195.5
$ws.Range("H139").Value = 2622
$ws.Range("I139").Value = 1996.3334
$ws.Range("K139").Value = 5989.0002
$ws.Range("M139").Value = -849.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 402.27274
$ws.Range("I2").Value = 112.666664
$ws.Range("J2").Value = 749.8
$ws.Range("K2").Value = 112.666664
$ws.Range("L2").Value = 749.8
$ws.Range("M2").Value = 0.3333360000000027
$ws.Range("N2").Value = -975.8
$ws.Range("H113").Value = 5191.75
$ws.Range("I113").Value = 5097
$ws.Range("K113").Value = 5097
$ws.Range("M113").Value = -2927
$ws.Range("H126").Value = 4541
$ws.Range("J126").Value = 4816.143
$ws.Range("L126").Value = 14448.429
$ws.Range("N126").Value = -19388.429
$ws.Range("H132").Value = 2024.625
$ws.Range("I132").Value = 1102.091
$ws.Range("J132").Value = 4054.2
$ws.Range("K132").Value = 3306.273
$ws.Range("L132").Value = 12162.6
$ws.Range("M132").Value = -776.2729999999997
$ws.Range("N132").Value = -17222.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1470.2727
$ws.Range("I16").Value = 1584.5
$ws.Range("J16").Value = 328
$ws.Range("K16").Value = 1584.5
$ws.Range("L16").Value = 328
$ws.Range("M16").Value = -1414.5
$ws.Range("N16").Value = -668
$ws.Range("H74").Value = 20000
$ws.Range("I74").Value = 20000
$ws.Range("K74").Value = 20000
$ws.Range("M74").Value = -19002
$ws.Range("H77").Value = 20000
$ws.Range("I77").Value = 20000
$ws.Range("K77").Value = 60000
$ws.Range("M77").Value = -55008
$ws.Range("H93").Value = 1902.6111
$ws.Range("I93").Value = 1935.1177
$ws.Range("J93").Value = 1350
$ws.Range("K93").Value = 1935.1177
$ws.Range("L93").Value = 1350
$ws.Range("M93").Value = -687.1177
$ws.Range("N93").Value = -3846
$ws.Range("H132").Value = 1666.6666
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -470
$ws.Range("N132").Value = -11060
$ws.Range("H136").Value = 3999.25
$ws.Range("I136").Value = 2749.5
$ws.Range("J136").Value = 5249
$ws.Range("K136").Value = 8248.5
$ws.Range("L136").Value = 15747
$ws.Range("M136").Value = -5698.5
$ws.Range("N136").Value = -20847

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 50000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51872
$ws.Range("H78").Value = 50000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -159360
$ws.Range("H132").Value = 2088.25
$ws.Range("I132").Value = 1884.3334
$ws.Range("J132").Value = 2700
$ws.Range("K132").Value = 5653.0002
$ws.Range("L132").Value = 8100
$ws.Range("M132").Value = -3123.0002
$ws.Range("N132").Value = -13160
$ws.Range("H136").Value = 3435
$ws.Range("I136").Value = 3163.087
$ws.Range("J136").Value = 4998.5
$ws.Range("K136").Value = 9489.261
$ws.Range("L136").Value = 14995.5
$ws.Range("M136").Value = -6939.261
$ws.Range("N136").Value = -20095.5
